$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for several rows based on repulled data
$ws.Range("F2").Value = -6
$ws.Range("F3").Value = -4
$ws.Range("F4").Value = 3
$ws.Range("F6").Value = -5
